$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "42.880.76"
$ws.Cells.Item(2, 5).Value = "  +0.11%  "

$ws.Cells.Item(3, 4).Value = "2.534.89"
$ws.Cells.Item(3, 5).Value = "  -1.48%  "

$ws.Cells.Item(4, 4).NumberFormat = "@"
$ws.Cells.Item(4, 4).Value = "0.998"
$ws.Cells.Item(4, 5).Value = "  -0.12%  "

$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "305.53"
$ws.Cells.Item(5, 5).Value = "  +1.08%  "

$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "98.95"
$ws.Cells.Item(6, 5).Value = "  +6.71%  "

$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = "0.585"
$ws.Cells.Item(7, 5).Value = "  +1.90%  "

$ws.Cells.Item(8, 5).Value = "  +0.04%  "

$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = "0.546"
$ws.Cells.Item(9, 5).Value = "  +0.14%  "

$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = "37.04"
$ws.Cells.Item(10, 5).Value = "  +2.83%  "

$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = "0.0815"
$ws.Cells.Item(11, 5).Value = "  +0.54%  "

$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = "7.74"
$ws.Cells.Item(12, 5).Value = "  -0.20%  "

$ws.Cells.Item(13, 5).Value = "  -1.10%  "

$ws.Cells.Item(14, 4).Value = "2.920.16"
$ws.Cells.Item(14, 5).Value = "  -1.51%  "

$ws.Cells.Item(15, 4).Value = "2.541.75"
$ws.Cells.Item(15, 5).Value = "  -1.91%  "

$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = "15.11"
$ws.Cells.Item(16, 5).Value = "  +6.15%  "

$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = "0.869"
$ws.Cells.Item(17, 5).Value = "  -1.63%  "

$ws.Cells.Item(18, 4).Value = "42.914.01"

$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = "13.03"
$ws.Cells.Item(19, 5).Value = "  +2.59%  "

$ws.Cells.Item(20, 4).Value = "0.0₃0983"
$ws.Cells.Item(20, 5).Value = "  -0.87%  "

$ws.Cells.Item(21, 5).Value = "  -1.94%  "

$ws.Cells.Item(22, 5).Value = "  -0.40%  "

$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = "253.33"
$ws.Cells.Item(23, 5).Value = "  -0.08%  "

$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = "3.13"
$ws.Cells.Item(24, 5).Value = "  +6.18%  "

$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = "2.05"
$ws.Cells.Item(25, 5).Value = "  -3.59%  "

$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = "26.91"
$ws.Cells.Item(26, 5).Value = "  -6.93%  "

$ws.Cells.Item(27, 5).Value = "  +0.20%  "

$ws.Cells.Item(28, 2).Value = "Toncoin"
$ws.Cells.Item(28, 3).Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = "2.34"
$ws.Cells.Item(28, 5).Value = "  +10.19%  "

$ws.Cells.Item(29, 2).Value = "Cosmos"
$ws.Cells.Item(29, 3).Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = "10.45"
$ws.Cells.Item(29, 5).Value = "  +1.70%  "

$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = "38.43"
$ws.Cells.Item(30, 5).Value = "  +4.03%  "

$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = "6.11"
$ws.Cells.Item(31, 5).Value = "  +1.63%  "

$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = "158.10"
$ws.Cells.Item(32, 5).Value = "  +2.39%  "

$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = "3.33"
$ws.Cells.Item(33, 5).Value = "  -1.97%  "

$ws.Cells.Item(34, 5).Value = "  -2.70%  "

$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = "0.0795"
$ws.Cells.Item(35, 5).Value = "  -0.51%  "

$ws.Cells.Item(36, 5).Value = "  -4.12%  "

$ws.Cells.Item(37, 2).Value = "Kaspa"
$ws.Cells.Item(37, 3).Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = "0.116"
$ws.Cells.Item(37, 5).Value = "  +2.06%  "

$ws.Cells.Item(38, 2).Value = "Celestia"
$ws.Cells.Item(38, 3).Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = "18.26"
$ws.Cells.Item(38, 5).Value = "  -0.28%  "

$ws.Cells.Item(39, 5).Value = "  +0.54%  "

$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = "24.12"
$ws.Cells.Item(40, 5).Value = "  +2.96%  "

$ws.Cells.Item(41, 2).Value = "NEARProtocol"
$ws.Cells.Item(41, 3).Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = "3.44"
$ws.Cells.Item(41, 5).Value = "  +1.08%  "

$ws.Cells.Item(42, 2).Value = "RenderToken"
$ws.Cells.Item(42, 3).Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = "3.91"
$ws.Cells.Item(42, 5).Value = "  +1.01%  "

$ws.Cells.Item(43, 2).Value = "ApeXProtocol"
$ws.Cells.Item(43, 3).Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = "2.09"
$ws.Cells.Item(43, 5).Value = "  +1.33%  "

$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = "0.0304"
$ws.Cells.Item(44, 5).Value = "  -2.55%  "

$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = "0.997"
$ws.Cells.Item(45, 5).Value = "  -0.03%  "

$ws.Cells.Item(46, 4).Value = "2.045.11"
$ws.Cells.Item(46, 5).Value = "  -1.97%  "

$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = "86.16"
$ws.Cells.Item(47, 5).Value = "  +1.56%  "

$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = "8.97"
$ws.Cells.Item(48, 5).Value = "  -3.25%  "

$ws.Cells.Item(49, 4).Value = "2.779.38"
$ws.Cells.Item(49, 5).Value = "  -1.35%  "

$ws.Cells.Item(50, 2).Value = "Algorand"
$ws.Cells.Item(50, 3).Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = "0.191"
$ws.Cells.Item(50, 5).Value = "  +0.26%  "

$ws.Cells.Item(51, 2).Value = "Aave"
$ws.Cells.Item(51, 3).Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = "102.94"
$ws.Cells.Item(51, 5).Value = "  -4.00%  "

